# Applies the 'Updated cryptos list' data refresh to Sheet1.
# Each target cell originally holds an inline text string (t="inlineStr"),
# including cells whose text happens to look numeric (e.g. "328.04", "1.000").
# A plain  Range.Value = '...'  assignment lets Excel's COM layer auto-coerce
# such strings into real numbers, which would silently change the cell's type.
# Forcing the NumberFormat to text ("@") before the write keeps the value a
# literal string; resetting the Style to "Normal" afterwards drops the style
# index that NumberFormat introduced, so no cell picks up formatting it didn't
# have before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '27.927.72'
Set-TextValue 'E2' '  +0.85%  '
Set-TextValue 'D3' '1.761.99'
Set-TextValue 'E3' '  -0.71%  '
Set-TextValue 'E4' '  -0.09%  '
Set-TextValue 'D5' '328.25'
Set-TextValue 'E5' '  +0.66%  '
Set-TextValue 'D6' '0.9999'
Set-TextValue 'E6' '  -0.04%  '
Set-TextValue 'D7' '0.4648'
Set-TextValue 'E7' '  +1.22%  '
Set-TextValue 'D8' '0.3513'
Set-TextValue 'E8' '  -2.13%  '
Set-TextValue 'D9' '43.74'
Set-TextValue 'E9' '  +4.21%  '
Set-TextValue 'D10' '0.07360'
Set-TextValue 'E10' '  -1.59%  '
Set-TextValue 'D11' '1.083'
Set-TextValue 'E11' '  -1.72%  '
Set-TextValue 'D12' '1.000'
Set-TextValue 'E12' '  -0.13%  '
Set-TextValue 'D13' '20.56'
Set-TextValue 'E13' '  -1.25%  '
Set-TextValue 'D14' '5.994'
Set-TextValue 'E14' '  -0.68%  '
Set-TextValue 'E15' '  -0.91%  '
Set-TextValue 'D16' '1.764.29'
Set-TextValue 'E16' '  -0.74%  '
Set-TextValue 'D17' '92.45'
Set-TextValue 'E17' '  -1.23%  '
Set-TextValue 'E18' '  -0.44%  '
Set-TextValue 'E19' '  +0.21%  '
Set-TextValue 'D20' '0.9997'
Set-TextValue 'E20' '  -0.01%  '
Set-TextValue 'E21' '  -1.37%  '
Set-TextValue 'D22' '5.757'
Set-TextValue 'E22' '  -0.56%  '
Set-TextValue 'D23' '27.948.05'
Set-TextValue 'D24' '11.12'
Set-TextValue 'E24' '  -1.42%  '
Set-TextValue 'D25' '2.154'
Set-TextValue 'E25' '  +3.34%  '
Set-TextValue 'D26' '162.34'
Set-TextValue 'E26' '  -1.39%  '
Set-TextValue 'E27' '  -1.74%  '
Set-TextValue 'D28' '1.965.24'
Set-TextValue 'E28' '  -0.87%  '
Set-TextValue 'D29' '2.171'
Set-TextValue 'E29' '  +0.38%  '
Set-TextValue 'D30' '123.20'
Set-TextValue 'E30' '  -2.37%  '
Set-TextValue 'D31' '1.066'
Set-TextValue 'E31' '  -2.88%  '
Set-TextValue 'D32' '0.09286'
Set-TextValue 'E32' '  +0.75%  '
Set-TextValue 'D33' '3.643'
Set-TextValue 'E33' '  -0.88%  '
Set-TextValue 'D34' '5.548'
Set-TextValue 'E34' '  +0.31%  '
Set-TextValue 'B35' 'Aptos'
Set-TextValue 'C35' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D35' '11.65'
Set-TextValue 'E35' '  -1.23%  '
Set-TextValue 'B36' 'VeChain'
Set-TextValue 'C36' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D36' '0.02268'
Set-TextValue 'E36' '  -1.07%  '
Set-TextValue 'D37' '0.06069'
Set-TextValue 'E37' '  -0.63%  '
Set-TextValue 'D38' '0.2060'
Set-TextValue 'E38' '  -1.35%  '
Set-TextValue 'D39' '4.901'
Set-TextValue 'E39' '  -1.20%  '
Set-TextValue 'D40' '0.6120'
Set-TextValue 'E40' '  -2.99%  '
Set-TextValue 'D41' '1.179'
Set-TextValue 'E41' '  -0.02%  '
Set-TextValue 'D42' '7.771'
Set-TextValue 'E42' '  +0.07%  '
Set-TextValue 'D43' '1.354'
Set-TextValue 'E43' '  -2.93%  '
Set-TextValue 'D44' '13.09'
Set-TextValue 'E44' '  -0.62%  '
Set-TextValue 'D45' '3.735'
Set-TextValue 'E45' '  +0.18%  '
Set-TextValue 'D46' '0.5782'
Set-TextValue 'E46' '  -1.71%  '
Set-TextValue 'D47' '122.61'
Set-TextValue 'E47' '  +0.26%  '
Set-TextValue 'E48' '  -1.12%  '
Set-TextValue 'D49' '0.06811'
Set-TextValue 'E49' '  -1.96%  '
Set-TextValue 'D50' '1.124'
Set-TextValue 'E50' '  -1.23%  '
Set-TextValue 'D51' '72.04'
Set-TextValue 'E51' '  -0.33%  '
